# Weekly driver report update for 2025-04-20
# The "Good Drivers" table (rows 13-18) is re-sorted by Driver Vintage
# (most-recent-first) and the sample counts / good-roaming % / vintage
# dates are refreshed for this week's pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final contents for A13:E18, top to bottom (sorted by vintage, newest first).
$rows = @(
    @{ Row = 13; Adapter = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4";  Samples = 445055; Pct = 99.90000000000001; Vintage = "2024-11-10" },
    @{ Row = 14; Adapter = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9";   Samples = 77849;  Pct = 99.90000000000001; Vintage = "2021-08-18" },
    @{ Row = 15; Adapter = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1";   Samples = 34244;  Pct = 100;                Vintage = "2021-04-27" },
    @{ Row = 16; Adapter = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2";  Samples = 59673;  Pct = 100;                Vintage = "2020-08-05" },
    @{ Row = 17; Adapter = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6";   Samples = 113652; Pct = 100;                Vintage = "2020-01-06" },
    @{ Row = 18; Adapter = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1";   Samples = 56018;  Pct = 100;                Vintage = "2019-12-14" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.Adapter
    $ws.Cells.Item($rowNum, 2).Value = $r.Samples
    $ws.Cells.Item($rowNum, 4).Value = $r.Pct

    # Driver Vintage is stored as a literal text date ("2024-11-10"), not an
    # Excel date serial, so force the cell to Text before writing it.
    $vintageCell = $ws.Cells.Item($rowNum, 5)
    $vintageCell.NumberFormat = "@"
    $vintageCell.Value = $r.Vintage
}
